# Move the "Inclusion comments" values from K4/K5 to K6/K7
# (the comments belonged on the rows below where they were originally placed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing comment text before clearing the source cells.
$k4Value = $ws.Range("K4").Text
$k5Value = $ws.Range("K5").Text

# Clear the old (incorrectly placed) cells.
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()

# Write the comments onto the correct rows.
$ws.Range("K6").Value = $k4Value
$ws.Range("K7").Value = $k5Value

# Update the active selection to match the saved view state.
$ws.Range("J9").Select()
